$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (D = Fecha serial, J = Volumen, K = Precio minimo,
# L = Precio maximo, M = Precio promedio ponderado, P = Precio $/Kg)
$data = @{
    2  = @(44377, 40,  14000, 15000, 14500, 1115)
    3  = @(44350, 40,  23000, 25000, 24000, 1846)
    4  = @(44327, 50,  24000, 25000, 24400, 1877)
    5  = @(44383, 50,  15000, 16000, 15400, 1185)
    6  = @(44425, 60,  14000, 15000, 14500, 1115)
    7  = @(44308, 50,  26000, 27000, 26400, 2031)
    8  = @(44355, 60,  18000, 20000, 19000, 1462)
    9  = @(44159, 60,  30000, 32000, 31000, 2385)
    10 = @(44320, 50,  26000, 28000, 26800, 2062)
    11 = @(44334, 50,  26000, 28000, 27200, 2092)
    12 = @(44362, 40,  15000, 16000, 15500, 1192)
    13 = @(44433, 100, 13000, 14000, 13500, 1038)
    14 = @(44264, 40,  30000, 32000, 31000, 2385)
    15 = @(44316, 50,  27000, 28000, 27400, 2108)
    16 = @(44313, 50,  25000, 26000, 25600, 1969)
    17 = @(44435, 100, 13000, 14000, 13500, 1038)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value  = $vals[0]   # D - Fecha
    $ws.Cells.Item($row, 10).Value = $vals[1]   # J - Volumen
    $ws.Cells.Item($row, 11).Value = $vals[2]   # K - Precio minimo
    $ws.Cells.Item($row, 12).Value = $vals[3]   # L - Precio maximo
    $ws.Cells.Item($row, 13).Value = $vals[4]   # M - Precio promedio ponderado
    $ws.Cells.Item($row, 16).Value = $vals[5]   # P - Precio $/Kg
}
